$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 103 (shifts existing rows 103:202 down to 104:203)
$ws.Rows.Item(103).Insert()

# Populate the new row 103 with the new price-record data, matching the
# surrounding rows' constant columns (A,B,C,E,F,G,H,I,N,O,Q,R) plus the
# new record's own date/volume/price data.
$ws.Range("A103").Value = 4
$ws.Range("B103").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C103").Value = "Los Lagos"
$ws.Range("D103").Value = 44566
$ws.Range("D103").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E103").Value = 10
$ws.Range("F103").Value = 100112043
$ws.Range("G103").Value = "Pepino ensalada"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 70
$ws.Range("K103").Value = 10000
$ws.Range("L103").Value = 10000
$ws.Range("M103").Value = 10000
$ws.Range("N103").Value = "$/caja 60 unidades"
$ws.Range("O103").Value = "Región de Arica y Parinacota"
$ws.Range("P103").Value = 167
$ws.Range("Q103").Value = 60
$ws.Range("R103").Value = "Hortaliza"
